$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.459612070389937
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 0.1575252929769615
$ws.Range("E2").Value = 8.660232485948974
$ws.Range("G2").Value = 11.945164432584

$ws.Range("B3").Value = 1.459612070389937
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 337.1190423067083
$ws.Range("E3").Value = 8.660232485948974
$ws.Range("G3").Value = 348.9066814463153
